$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by every Status cell across all three sheets, so
#    every cell currently holding it needs to be rewritten with the new text.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$ws1.Range("B2").Value = $newStatus
$ws1.Range("C2").Value = $newStatus
$ws1.Range("B3").Value = $newStatus
$ws1.Range("C3").Value = $newStatus

$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) zh-cn sheet ("Handed back" in sync with en-US for zh-cn at 02:51:27) -
#    Latest Handback DateTime (column H) updates for both data rows.
# ---------------------------------------------------------------------------
$ws2.Range("H2").Value = "2016-03-21 02:51:27"
$ws2.Range("H3").Value = "2016-03-21 02:51:27"

# ---------------------------------------------------------------------------
# 3) de-de sheet - Latest Handback DateTime (column H) updates for both rows,
#    using a distinct timestamp from the zh-cn handback.
# ---------------------------------------------------------------------------
$ws3.Range("H2").Value = "2016-03-21 02:51:33"
$ws3.Range("H3").Value = "2016-03-21 02:51:33"

# ---------------------------------------------------------------------------
# 4) New columns F (Latest Target File) / G (Latest Handback File) populated
#    with hyperlinks, for both the zh-cn and de-de sheets, rows 2 and 3.
#    F repeats the source markdown file link (same as column A); G repeats
#    the localized target file link (same as column D).
# ---------------------------------------------------------------------------

# zh-cn, row 2 (92a7e3e7...)
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/eeb3416dbba13d77c93a0e3e977696210c78f629/e2e/92a7e3e7-d8df-402f-9014-050ef4a47fe6.md", "", "", "92a7e3e7-d8df-402f-9014-050ef4a47fe6.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6a62964f00e5d629cdd76bd7cdff611b0a8fc45c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/92a7e3e7-d8df-402f-9014-050ef4a47fe6.4d2c580fbfb8e49e7ef9da01482b087221e5382b.zh-cn.xlf", "", "", "92a7e3e7-d8df-402f-9014-050ef4a47fe6.4d2c580fbfb8e49e7ef9da01482b087221e5382b.zh-cn.xlf")

# zh-cn, row 3 (caca818b...)
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/eeb3416dbba13d77c93a0e3e977696210c78f629/e2e/caca818b-7dc9-44b7-a836-98706a65cf53.md", "", "", "caca818b-7dc9-44b7-a836-98706a65cf53.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6a62964f00e5d629cdd76bd7cdff611b0a8fc45c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/caca818b-7dc9-44b7-a836-98706a65cf53.9d3c212640697a9ba3122ff25d55785288c989b0.zh-cn.xlf", "", "", "caca818b-7dc9-44b7-a836-98706a65cf53.9d3c212640697a9ba3122ff25d55785288c989b0.zh-cn.xlf")

# de-de, row 2 (92a7e3e7...)
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/eeb3416dbba13d77c93a0e3e977696210c78f629/e2e/92a7e3e7-d8df-402f-9014-050ef4a47fe6.md", "", "", "92a7e3e7-d8df-402f-9014-050ef4a47fe6.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/48411b40cddbeff54237c504a304c9fb58def7b9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/92a7e3e7-d8df-402f-9014-050ef4a47fe6.4d2c580fbfb8e49e7ef9da01482b087221e5382b.de-de.xlf", "", "", "92a7e3e7-d8df-402f-9014-050ef4a47fe6.4d2c580fbfb8e49e7ef9da01482b087221e5382b.de-de.xlf")

# de-de, row 3 (caca818b...)
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/eeb3416dbba13d77c93a0e3e977696210c78f629/e2e/caca818b-7dc9-44b7-a836-98706a65cf53.md", "", "", "caca818b-7dc9-44b7-a836-98706a65cf53.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/48411b40cddbeff54237c504a304c9fb58def7b9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/caca818b-7dc9-44b7-a836-98706a65cf53.9d3c212640697a9ba3122ff25d55785288c989b0.de-de.xlf", "", "", "caca818b-7dc9-44b7-a836-98706a65cf53.9d3c212640697a9ba3122ff25d55785288c989b0.de-de.xlf")
